# Apply the "Rebuild TOC and Index" / "Rebuild PDF" reorg edit.
#
# Before:
#   Para A (ilvl=2): "Rebuild PDF"
#   Para B (ilvl=0): "Copy all referenced PDFs to " + [CodeChar]"Users\Public\Documentation"
#
# After:
#   Para A (ilvl=2): "Rebuild TOC and Index"
#   Para B (ilvl=2): "Rebuild PDF"   (CodeChar run removed entirely)

$d = $word.ActiveDocument

# --- Paragraph A: the list item that currently reads exactly "Rebuild PDF" ---
$rngA = $d.Content
$rngA.Find.Execute("Rebuild PDF", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rngA.Text = "Rebuild TOC and Index"

# --- Paragraph B: the list item "Copy all referenced PDFs to Users\Public\Documentation" ---
# Replacing the whole range text collapses both runs (incl. the CodeChar-styled one) into a
# single plain run, matching the target content.
$rngB = $d.Content
$rngB.Find.Execute("Copy all referenced PDFs to Users\Public\Documentation", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rngB.Text = "Rebuild PDF"

# --- Promote paragraph B to list level 3 (ilvl=2, 0-based) to match paragraph A ---
# At this point paragraph A already reads "Rebuild TOC and Index", so "Rebuild PDF" is unique.
$rngC = $d.Content
$rngC.Find.Execute("Rebuild PDF", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rngC.ListFormat.ListLevelNumber = 3
